# "update segmentation with new timing"
# LWP2_0009_lab_timing.xlsx - "In Lab" sheet
#
# The edit re-times several lab-session segmentation events in column B
# (now carrying second-level precision, so the displayed format moves from
# h:mm to h:mm:ss), and moves the saved selection/scroll position.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("In Lab")

# --- Updated segmentation timestamps (column B), with second precision ---
$ws.Range("B7").Value2  = 0.63368055555555558
$ws.Range("B9").Value2  = 0.64188657407407412
$ws.Range("B11").Value2 = 0.65130787037037041
$ws.Range("B12").Value2 = 0.65216435185185184
$ws.Range("B13").Value2 = 0.65305555555555561
$ws.Range("B14").Value2 = 0.66052083333333333
$ws.Range("B15").Value2 = 0.66145833333333337
$ws.Range("B16").Value2 = 0.66315972222222219
$ws.Range("B17").Value2 = 0.66548611111111111
$ws.Range("B18").Value2 = 0.67282407407407396
$ws.Range("B20").Value2 = 0.67391203703703706
$ws.Range("B21").Value2 = 0.67651620370370369
$ws.Range("B22").Value2 = 0.67800925925925926
$ws.Range("B23").Value2 = 0.67871527777777774
$ws.Range("B24").Value2 = 0.67913194444444447

# Every timing cell in the column (including the ones whose value did not
# move, e.g. B6/B8/B10/B19/B25) now renders with seconds.
$ws.Range("B6:B25").NumberFormat = "h:mm:ss"

# --- Saved view state: selection moves to B25 (drops the old topLeftCell="C1" scroll) ---
$ws.Range("B25").Select() | Out-Null
